# Add a new entry (row 24) to the "Activités" journal table:
# Date 2021-03-12, Début 13:36, Fin 15:03, Activité "Rédaction documentation",
# Commentaires "Avancement sur la documentation".
# Also move the active-cell selection to G25 (one row below the new entry).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activités")

# A24: Date (Excel serial 44267 = 2021-03-12)
$ws.Cells.Item(24, 1).Value = 44267

# B24: Début (time, Excel serial fraction = 13:36)
$ws.Cells.Item(24, 2).Value = 0.56666666666666665

# C24: Fin (time, Excel serial fraction = 15:03)
$ws.Cells.Item(24, 3).Value = 0.62708333333333333

# E24: Activité
$ws.Cells.Item(24, 5).Value = "Rédaction documentation"

# G24: Commentaires
$ws.Cells.Item(24, 7).Value = "Avancement sur la documentation"

# Move selection to G25, matching the diff's sheetView change.
$ws.Range("G25").Select()

$excel.ActiveWorkbook.Save()
